$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new worksheet "IESD-AAaWER" at the end of the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "IESD-AAaWER"

# ---------------------------------------------------------------------------
# 2. Update the "About" sheet: insert a new row (3) with the new sheet's
#    title, and append a new section at the bottom describing the new
#    "Average Age at Which Equipment Retires" output.
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Insert a blank row at row 3 (shifts everything below down by one), then
# populate it with the title for the new sheet (matches formatting of the
# other sheet-title rows above it).
$about.Rows.Item(3).Insert()
$about.Range("A3").Value = "IESD-AAaWER Average Age at Which Equipment Retires"

# New explanatory section appended at the bottom of the sheet.
$about.Range("A27").Value = "Average Age at Which Equipment Retires"
$about.Range("A27").Font.Bold = $true
$about.Range("A27:B27").Interior.Color = 0xECCAA6

$about.Range("A28").Value = "The output ""Average Age at Which Equipment Retires"" is used to help levelize capital costs across"
$about.Range("A29").Value = "units of energy consumed by equipment over the lifetime of that equipment when firms make"
$about.Range("A30").Value = "decisions about the type of fuel newly purchased equipment should used."
$about.Range("A31").Value = "It is not used directly in equipment tracking (which is based on the full survival curve,"
$about.Range("A32").Value = "not an average value)."

# ---------------------------------------------------------------------------
# 3. Populate the new "IESD-AAaWER" sheet.
# ---------------------------------------------------------------------------
$newSheet.Tab.Color = 0xECCAA6

$newSheet.Columns.Item(1).ColumnWidth = 23.21875
$newSheet.Columns.Item(2).ColumnWidth = 24.21875

$newSheet.Range("A1").Value = "Unit: years"
$newSheet.Range("A1").Font.Italic = $true

$newSheet.Range("A2").Value = "Industrial equipment"

$newSheet.Range("B1").Value = "Average age at retirement"
$newSheet.Range("B1").HorizontalAlignment = -4152

$newSheet.Range("B2").Formula = "=XLOOKUP(0.5,Calculations!B5:AZ5,Calculations!B1:AZ1,,1)"
